$d = $word.ActiveDocument

# "Time limit: 200 ms." -> "Time limit: 300 ms."
$d.Content.Find.Execute("200", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "300", 2)
